$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Temporarily un-bold the header row so table creation does not capture
# the existing bold formatting as a header dxf override, then restore it.
$ws.Range("A1:U1").Font.Bold = $false
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U66"), $null, 1)
$ws.Range("A1:U1").Font.Bold = $true
Write-Host "done"
